$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell K1, matching the style of the existing header row (J1)
$ws.Range("K1").Value = "PhylogenySorting"
$ws.Range("K1").Font.Bold = $true
$ws.Range("K1").HorizontalAlignment = -4108

# Add new data cell K3, matching value of neighboring cell J3 ("Unassigned")
$ws.Range("K3").Value = "Unassigned"
